$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Account Number becomes a real number (was text), Netpay amount updated
$ws.Range("A2").Value = 32145698741
$ws.Range("C2").Value = 40989.1

# New row 3: another employee's leave/pay record
$ws.Range("A3").Value = 123654789963
$ws.Range("B3").Value = "Vidya Sagar pogiri"
$ws.Range("C3").Value = 9793.33
$ws.Range("D3").Value = "November"

# Column A widened to fit the new, wider numeric account numbers
$ws.Range("A1").EntireColumn.ColumnWidth = 15.6
